$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: a new "2022-Q4" row is inserted at the
#    top of the data (row 2), pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Row 7 is brand new (the sheet previously only had 6 rows) - extend the
# index-column style from A6 down to it before filling in values.
$summary.Range("A6").Copy()
$summary.Range("A7").PasteSpecial(-4122)

$summaryData = @(
  @(0, "2022-Q4", 7, 1.07),
  @(1, "2022-Q3", 1, 0.01),
  @(2, "2022-Q2", 3, 0.04),
  @(3, "2022-Q1", 2, 0.32),
  @(4, "2021-Q4", 3, 0.54),
  @(5, "2021-Q3", 3, 0.58)
)

$r = 2
foreach ($row in $summaryData) {
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = "'" + $row[1]
    $summary.Cells.Item($r, 2).Style = "Normal"
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q4" fund-holdings sheet right after "总计" (i.e.
#    before the current "2022-Q3" sheet), reusing the "2022-Q3" sheet's
#    layout/formatting as a template, then overwrite its contents.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q3")
$template.Copy($template)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Extend the style used on A2 (index column) down through A3:A8.
$q4.Range("A2").Copy()
$q4.Range("A3:A8").PasteSpecial(-4122)

$q4Data = @(
  @(0, "005106", "银华农业产业股票A", "12.40", "93.03", "5.29", "0.6560", 8),
  @(1, "013142", "华商乐享互联灵活配置混合C", "5.21", "88.52", "3.13", "0.1631", 6),
  @(2, "001959", "华商乐享互联灵活配置混合A", "4.50", "88.52", "3.13", "0.1408", 6),
  @(3, "014064", "银华农业产业股票C", "1.12", "93.03", "5.29", "0.0592", 8),
  @(4, "001219", "上投摩根动态多因子策略混合A", "0.97", "92.08", "4.15", "0.0403", 2),
  @(5, "012430", "农银汇理瑞康6个月持有期混合", "1.16", "28.06", "1.11", "0.0129", 5),
  @(6, "017176", "上投摩根动态多因子策略混合C", "0.00", "92.08", "4.15", $null, 2)
)

$r = 2
foreach ($row in $q4Data) {
    $q4.Cells.Item($r, 1).Value = $row[0]

    $q4.Cells.Item($r, 2).Value = "'" + $row[1]
    $q4.Cells.Item($r, 2).Style = "Normal"

    $q4.Cells.Item($r, 3).Value = $row[2]

    $q4.Cells.Item($r, 4).Value = "'" + $row[3]
    $q4.Cells.Item($r, 4).Style = "Normal"

    $q4.Cells.Item($r, 5).Value = "'" + $row[4]
    $q4.Cells.Item($r, 5).Style = "Normal"

    $q4.Cells.Item($r, 6).Value = "'" + $row[5]
    $q4.Cells.Item($r, 6).Style = "Normal"

    if ($row[6] -eq $null) {
        $q4.Cells.Item($r, 7).Value = 0
    } else {
        $q4.Cells.Item($r, 7).Value = "'" + $row[6]
        $q4.Cells.Item($r, 7).Style = "Normal"
    }

    $q4.Cells.Item($r, 8).Value = $row[7]

    $r++
}
